# Change the lab date from 3/9/2021 to 3/11/2021 and move Word's "last
# edit" bookmark (_GoBack) to sit right after the new date fragment,
# exactly where the cursor would be after typing "11" over "9".
#
# Note: this engine re-coalesces adjacent runs that share identical
# formatting whenever a paragraph's text is edited, *except* across a
# bookmark boundary. To keep "Lab Exerc" from re-merging with the run
# that holds the date while we edit it, we temporarily bookmark that
# boundary, make the edit, then remove the temporary bookmark (removing
# a bookmark is a structural change and does not trigger a run merge).

$d = $word.ActiveDocument

# 1) Protect the "Lab Exerc" / "ise 3/9/2021" run boundary with a throwaway bookmark.
$findLabExerc = $d.Content
$findLabExerc.Find.Execute("Lab Exerc") | Out-Null
$d.Bookmarks.Add("zzTemp", $d.Range($findLabExerc.End, $findLabExerc.End))

# 2) Drop the real _GoBack bookmark where the cursor lands after retyping
#    the day: right after "ise 3/9", before "/2021". Adding a bookmark
#    with a name that already exists elsewhere relocates it, so this also
#    clears the old _GoBack from the second paragraph.
$findDay = $d.Content
$findDay.Find.Execute("ise 3/9") | Out-Null
$d.Bookmarks.Add("_GoBack", $d.Range($findDay.End, $findDay.End))

# 3) Fix the day itself, 9 -> 11. This edit is sandwiched between the two
#    bookmarks, so it cannot merge into the runs on either side.
$findNine = $d.Content
$findNine.Find.Execute("3/9") | Out-Null
$d.Range($findNine.End - 1, $findNine.End).Text = "11"

# 4) Drop the scratch bookmark now that the edit is safely isolated.
$d.Bookmarks("zzTemp").Delete()
